$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 41
$ws.Cells.Item($row, 1).Value = "05/01/2026 04:43:18"
$ws.Cells.Item($row, 2).Value = "05/01 04:15"
$ws.Cells.Item($row, 3).Value = "Metrópoles"
$ws.Cells.Item($row, 4).Value = "Professor teria sido morto em emboscada por homem que conheceu em app"
$ws.Cells.Item($row, 5).Value = "https://www.metropoles.com/distrito-federal/na-mira/professor-teria-sido-morto-em-emboscada-armada-por-homem-que-conheceu-em-app"
$ws.Cells.Item($row, 6).Value = "lula"
$ws.Cells.Item($row, 7).Value = "Corpo de João Emmanuel Moura foi encontrado em parada de ônibus com o ce&lt;b&gt;lula&lt;/b&gt;r, que agora é periciado pela polícia"
